# Update "list of issues" sheet: append 5 new country rows that reuse the
# existing "Issues with getting the data from OSM. Cables again" text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$problemText = "Issues with getting the data from OSM. Cables again"
$countries = @("Chad", "Cameroon", "Central African Republic", "South Sudan", "Gabon")

$startRow = 9
for ($i = 0; $i -lt $countries.Length; $i++) {
    $r = $startRow + $i
    $ws.Range("A$r").Value = $countries[$i]
    $ws.Range("B$r").Value = $problemText

    # Match the wrap-text styling used by the rest of the table (style index 1).
    $ws.Range("A$r").WrapText = $true
    $ws.Range("B$r").WrapText = $true

    # Match the row height used by the other wrapped rows (6-8) above.
    $ws.Rows.Item($r).RowHeight = 48
}

# Move the active selection to the last new cell, like the source workbook.
$ws.Range("B13").Select() | Out-Null
